$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Update the Username (was "azar") and Password (was "demo123") test data values
$ws.Range("B2").Value = "nirai"
$ws.Range("B3").Value = "nirai123"

# Move the active selection on the Login sheet from F10 to C4
$ws.Activate()
$ws.Range("C4").Select()
